$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Scenario/Status/Browser/Timestamp log rows appended to the ScenarioStatus
# sheet (DB login and Login func fixed) - rows 280-313.
$newRows = @(
  @("Add Discounts Positive", "PASSED", "chrome", "11-11-2022 17:02:05"),
  @("Add Discount Negative", "PASSED", "chrome", "11-11-2022 17:02:33"),
  @("Edit Discounts", "PASSED", "chrome", "11-11-2022 17:03:13"),
  @("Delete Discounts Positive", "PASSED", "chrome", "11-11-2022 17:03:55"),
  @("Delete Discounts Negative", "PASSED", "chrome", "11-11-2022 17:04:37"),
  @("Add Nationality Positive", "PASSED", "chrome", "11-11-2022 17:04:59"),
  @("Add Nationality Negative", "PASSED", "chrome", "11-11-2022 17:05:19"),
  @("Edit Nationality", "PASSED", "chrome", "11-11-2022 17:05:50"),
  @("Delete Nationality Positive", "FAILED", "chrome", "11-11-2022 17:05:58"),
  @("Delete Nationality Negative", "FAILED", "chrome", "11-11-2022 17:07:02"),
  @("Add Discounts Positive", "FAILED", "chrome", "11-11-2022 17:16:34"),
  @("Add Discounts Positive", "FAILED", "chrome", "11-11-2022 17:17:56"),
  @("Add Discounts Positive", "FAILED", "chrome", "14-11-2022 08:18:56"),
  @("Add Discount Negative", "PASSED", "chrome", "14-11-2022 08:19:57"),
  @("Edit Discounts", "PASSED", "chrome", "14-11-2022 08:21:12"),
  @("Delete Discounts Positive", "PASSED", "chrome", "14-11-2022 08:22:12"),
  @("Delete Discounts Negative", "PASSED", "chrome", "14-11-2022 08:23:30"),
  @("Add Nationality Positive", "PASSED", "chrome", "14-11-2022 08:23:52"),
  @("Add Nationality Negative", "PASSED", "chrome", "14-11-2022 08:24:10"),
  @("Edit Nationality", "FAILED", "chrome", "14-11-2022 08:25:29"),
  @("Delete Nationality Positive", "FAILED", "chrome", "14-11-2022 08:26:43"),
  @("Delete Nationality Negative", "PASSED", "chrome", "14-11-2022 08:27:37"),
  @("Add Discounts Positive", "FAILED", "chrome", "14-11-2022 09:45:11"),
  @("Add Discount Negative", "PASSED", "chrome", "14-11-2022 09:45:32"),
  @("Edit Discounts", "FAILED", "chrome", "14-11-2022 09:45:42"),
  @("Delete Discounts Positive", "FAILED", "chrome", "14-11-2022 09:46:51"),
  @("Delete Discounts Negative", "PASSED", "chrome", "14-11-2022 09:47:29"),
  @("Add Nationality Positive", "PASSED", "chrome", "14-11-2022 09:47:46"),
  @("Add Nationality Negative", "PASSED", "chrome", "14-11-2022 09:48:04"),
  @("Edit Nationality", "PASSED", "chrome", "14-11-2022 09:48:36"),
  @("Delete Nationality Positive", "PASSED", "chrome", "14-11-2022 09:49:14"),
  @("Delete Nationality Negative", "PASSED", "chrome", "14-11-2022 09:49:38"),
  @("Nations testing with JDBC", "FAILED", "chrome", "15-11-2022 09:41:37"),
  @("Nations testing with JDBC", "PASSED", "chrome", "15-11-2022 09:43:10")

)

$startRow = 280
for ($i = 0; $i -lt $newRows.Length; $i++) {
  $r = $startRow + $i
  $row = $newRows[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
}
